$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: add columns P1 (14) and Q1 (15). Copy the format from the
# existing O1 header cell (bold, centered, thin border) so the new cells
# pick up the same cell style already used by the rest of row 1.
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: fix up columns I, K, M, O and append new columns P, Q.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2 (new)
}
